# GPLIM-3541: add Material Type as required header for Manifest uploads
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column G: "Material Type" - formatted like the other header cells (A1)
$ws.Range("G1").Value = "Material Type"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# New data column G2:G24: "DNA:Genomic" - formatted like the centered data cells (C2)
$ws.Range("G2:G24").Value = "DNA:Genomic"
$ws.Range("C2").Copy()
$ws.Range("G2:G24").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the active selection to the new column, matching the saved view state
[void]$ws.Activate()
[void]$ws.Range("G1:G24").Select()
